$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, C(new), D(new), E(new)
$updates = @(
    @(88,   71266, 12436,  110295157),
    @(91,   18857, 3383,   75163013),
    @(97,   22,    4,      354577),
    @(98,   6295,  911,    19452171),
    @(100,  9339,  1304,   23760393),
    @(121,  1306135, 220385, 2274587611),
    @(122,  366,   48,     1162215),
    @(129,  633363, 104966, 3427264164),
    @(132,  585656, 90780,  3462321824),
    @(136,  26677, 4272,   143561800),
    @(139,  76641, 17494,  114132538),
    @(144,  25074, 6170,   92430892),
    @(146,  7440,  1667,   37684204),
    @(151,  39928, 7155,   60363412),
    @(154,  18453, 3296,   72760118),
    @(156,  12398, 2144,   40047588),
    @(171,  95821, 18308,  490672924),
    @(189,  100465, 18353, 556064990),
    @(224,  39675, 7247,   260923034),
    @(229,  612545, 121244, 1040753491),
    @(240,  205901, 33984, 1068739496)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Range("C$row").Value = $u[1]
    $ws.Range("D$row").Value = $u[2]
    $ws.Range("E$row").Value = $u[3]
}
